$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.303.26"
$ws.Range("E2").Value = "'  +3.88%  "
$ws.Range("D3").Value = "'1.835.52"
$ws.Range("E3").Value = "'  +4.09%  "
$ws.Range("D4").Value = "'1.021"
$ws.Range("E4").Value = "'  +2.41%  "
$ws.Range("D5").Value = "'316.94"
$ws.Range("E5").Value = "'  +1.68%  "
$ws.Range("E6").Value = "'  +2.82%  "
$ws.Range("D7").Value = "'0.4339"
$ws.Range("E7").Value = "'  +1.83%  "
$ws.Range("D8").Value = "'0.3714"
$ws.Range("E8").Value = "'  +2.26%  "
$ws.Range("D9").Value = "'0.07305"
$ws.Range("E9").Value = "'  +2.13%  "
$ws.Range("D10").Value = "'0.8745"
$ws.Range("E10").Value = "'  +3.47%  "
$ws.Range("D11").Value = "'2.048.50"
$ws.Range("E11").Value = "'  +14.94%  "
$ws.Range("E12").Value = "'  +4.87%  "
$ws.Range("D13").Value = "'5.465"
$ws.Range("E13").Value = "'  +4.51%  "
$ws.Range("D14").Value = "'6.668"
$ws.Range("E14").Value = "'  +3.97%  "
$ws.Range("D15").Value = "'0.07131"
$ws.Range("E15").Value = "'  +4.04%  "
$ws.Range("D16").Value = "'82.02"
$ws.Range("E16").Value = "'  +4.30%  "
$ws.Range("D17").Value = "'1.026"
$ws.Range("E17").Value = "'  +3.04%  "
$ws.Range("D18").Value = "'0.000008992"
$ws.Range("E18").Value = "'  +3.92%  "
$ws.Range("D19").Value = "'1.017"
$ws.Range("E19").Value = "'  +2.18%  "
$ws.Range("D20").Value = "'15.41"
$ws.Range("E20").Value = "'  +2.94%  "
$ws.Range("D21").Value = "'27.351.95"
$ws.Range("E21").Value = "'  +4.00%  "
$ws.Range("D22").Value = "'5.229"
$ws.Range("E22").Value = "'  +2.89%  "
$ws.Range("E23").Value = "'  +0.23%  "
$ws.Range("D24").Value = "'2.271.64"
$ws.Range("E24").Value = "'  +13.54%  "
$ws.Range("D25").Value = "'156.40"
$ws.Range("E25").Value = "'  +3.58%  "
$ws.Range("D26").Value = "'1.893"
$ws.Range("E26").Value = "'  +1.25%  "
$ws.Range("D27").Value = "'18.53"
$ws.Range("E27").Value = "'  +2.87%  "
$ws.Range("D28").Value = "'5.280"
$ws.Range("E28").Value = "'  +3.92%  "
$ws.Range("D29").Value = "'1.926"
$ws.Range("E29").Value = "'  +7.93%  "
$ws.Range("D30").Value = "'115.35"
$ws.Range("E30").Value = "'  +1.62%  "
$ws.Range("D31").Value = "'0.09001"
$ws.Range("E31").Value = "'  +0.90%  "
$ws.Range("D32").Value = "'1.199"
$ws.Range("E32").Value = "'  +7.77%  "
$ws.Range("D33").Value = "'0.7581"
$ws.Range("E33").Value = "'  +4.41%  "
$ws.Range("D34").Value = "'4.458"
$ws.Range("E34").Value = "'  +3.53%  "
$ws.Range("D35").Value = "'2.838"
$ws.Range("E35").Value = "'  +4.30%  "
$ws.Range("E36").Value = "'  +2.60%  "
$ws.Range("D37").Value = "'1.145"
$ws.Range("E37").Value = "'  +6.69%  "
$ws.Range("D38").Value = "'0.01950"
$ws.Range("E38").Value = "'  +4.16%  "
$ws.Range("D39").Value = "'0.05256"
$ws.Range("E39").Value = "'  +2.33%  "
$ws.Range("E40").Value = "'  +5.34%  "
$ws.Range("D41").Value = "'2.784"
$ws.Range("E41").Value = "'  +8.96%  "
$ws.Range("D42").Value = "'0.1661"
$ws.Range("E42").Value = "'  +3.65%  "
$ws.Range("D43").Value = "'6.512"
$ws.Range("E43").Value = "'  +4.29%  "
$ws.Range("D44").Value = "'8.442"
$ws.Range("E44").Value = "'  +5.77%  "
$ws.Range("B45").Value = "'EnergySwap"
$ws.Range("C45").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'10.55"
$ws.Range("E45").Value = "'  +3.69%  "
$ws.Range("B46").Value = "'Quant"
$ws.Range("C46").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'107.72"
$ws.Range("E46").Value = "'  +3.14%  "
$ws.Range("E47").Value = "'  +3.38%  "
$ws.Range("D48").Value = "'0.4631"
$ws.Range("E48").Value = "'  +3.70%  "
$ws.Range("E49").Value = "'  +9.85%  "
$ws.Range("E50").Value = "'  +3.90%  "
$ws.Range("D51").Value = "'0.06273"
$ws.Range("E51").Value = "'  +1.47%  "
